$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.694.20'
$ws.Range("E2").Value = '  -2.36%  '

$ws.Range("D3").Value = '3.231.03'
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("E4").Value = '  +0.03%  '

$origStyle_D5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.94'
$ws.Range("D5").Style = $origStyle_D5
$ws.Range("E5").Value = '  -1.65%  '

$origStyle_D6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.01'
$ws.Range("D6").Style = $origStyle_D6
$ws.Range("E6").Value = '  -4.02%  '

$origStyle_D7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.624'
$ws.Range("D7").Style = $origStyle_D7
$ws.Range("E7").Value = '  -0.52%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").Value = '3.229.96'
$ws.Range("E9").Value = '  -1.37%  '

$ws.Range("E10").Value = '  -2.82%  '

$ws.Range("E11").Value = '  +0.33%  '

$ws.Range("E12").Value = '  -3.27%  '

$ws.Range("D13").Value = '3.792.01'
$ws.Range("E13").Value = '  -1.47%  '

$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("D15").Value = '64.803.86'
$ws.Range("E15").Value = '  -2.16%  '

$origStyle_D16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.72'
$ws.Range("D16").Style = $origStyle_D16
$ws.Range("E16").Value = '  -2.44%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.223.92'
$ws.Range("E17").Value = '  -1.85%  '

$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$origStyle_D18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000159'
$ws.Range("D18").Style = $origStyle_D18
$ws.Range("E18").Value = '  -3.03%  '

$origStyle_D19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '415.91'
$ws.Range("D19").Style = $origStyle_D19
$ws.Range("E19").Value = '  -3.92%  '

$origStyle_D20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.38'
$ws.Range("D20").Style = $origStyle_D20
$ws.Range("E20").Value = '  -2.56%  '

$origStyle_D21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.83'
$ws.Range("D21").Style = $origStyle_D21
$ws.Range("E21").Value = '  -2.47%  '

$origStyle_D22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.21'
$ws.Range("D22").Style = $origStyle_D22
$ws.Range("E22").Value = '  -2.39%  '

$ws.Range("E23").Value = '  -0.13%  '

$origStyle_D24 = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.25'
$ws.Range("D24").Style = $origStyle_D24
$ws.Range("E24").Value = '  -1.97%  '

$ws.Range("E25").Value = '  -0.66%  '

$ws.Range("E26").Value = '  +4.61%  '

$origStyle_D27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.496'
$ws.Range("D27").Style = $origStyle_D27
$ws.Range("E27").Value = '  -1.95%  '

$ws.Range("E28").Value = '  -2.04%  '

$origStyle_D29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.98'
$ws.Range("D29").Style = $origStyle_D29
$ws.Range("E29").Value = '  +1.67%  '

$ws.Range("E30").Value = '  +0.07%  '

$ws.Range("E31").Value = '  -5.03%  '

$origStyle_D32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.80'
$ws.Range("D32").Style = $origStyle_D32
$ws.Range("E32").Value = '  -2.06%  '

$origStyle_D33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = $origStyle_D33
$ws.Range("E33").Value = '  +0.06%  '

$ws.Range("E34").Value = '  -3.22%  '

$ws.Range("E35").Value = '  -2.93%  '

$ws.Range("E36").Value = '  -3.09%  '

$origStyle_D37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.55'
$ws.Range("D37").Style = $origStyle_D37
$ws.Range("E37").Value = '  -0.19%  '

$ws.Range("E38").Value = '  -1.69%  '

$ws.Range("D39").Value = '2.807.51'
$ws.Range("E39").Value = '  +0.98%  '

$ws.Range("E40").Value = '  -3.07%  '

$origStyle_D41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.48'
$ws.Range("D41").Style = $origStyle_D41
$ws.Range("E41").Value = '  -3.75%  '

$origStyle_D42 = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.20'
$ws.Range("D42").Style = $origStyle_D42
$ws.Range("E42").Value = '  -2.90%  '

$origStyle_D43 = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '39.39'
$ws.Range("D43").Style = $origStyle_D43
$ws.Range("E43").Value = '  -2.00%  '

$origStyle_D44 = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.723'
$ws.Range("D44").Style = $origStyle_D44
$ws.Range("E44").Value = '  -6.60%  '

$ws.Range("E45").Value = '  -4.35%  '

$origStyle_D46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0629'
$ws.Range("D46").Style = $origStyle_D46
$ws.Range("E46").Value = '  -4.56%  '

$ws.Range("E47").Value = '  -4.68%  '

$origStyle_D48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '301.91'
$ws.Range("D48").Style = $origStyle_D48
$ws.Range("E48").Value = '  -5.91%  '

$origStyle_D49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '21.99'
$ws.Range("D49").Style = $origStyle_D49
$ws.Range("E49").Value = '  -5.19%  '

$ws.Range("E50").Value = '  -1.46%  '
